$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells (AD1:AF1): copy the formatting from an existing header
# cell (bold, centered, bordered) so they match the rest of row 1, then
# set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-40: team record values (Wins=69, Losses=93, Ties=1) for every player.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 69   # column AD
    $ws.Cells.Item($r, 31).Value = 93   # column AE
    $ws.Cells.Item($r, 32).Value = 1    # column AF
}
